# "Minor update - GUI"
# Insert a new documentation row (pose_t_relative) into the hex_path sheet,
# right after the existing pose_t row, wrap the new NOTE cell, widen the
# columns that now hold long text, and move the active tab / selection from
# hex_obj to hex_path (matching the author's last on-screen state).

$wb = $excel.ActiveWorkbook

$wsObj  = $wb.Worksheets.Item("hex_obj")
$wsPath = $wb.Worksheets.Item("hex_path")

# --- hex_path: insert a new row 7 for "pose_t_relative" ------------------
$wsPath.Rows.Item(7).Insert()

$wsPath.Range("A7").Value = "pose_t_relative"
$wsPath.Range("B7").Value = "[6xN]"
$wsPath.Range("C7").Value = "end effector pose in world/base frame, relative to current pose."
$wsPath.Range("D7").Value = "NOTE: relative path is the default. E.g. All motion inputs will be interpreted as being relative to the current pose, unless the motion-planning dialog specifies 'absolute motion', in which case the pose_t array will be modified so that the motions are relative to the hexpod's home position."

# Taller row for the wrapped note, and wrap text on the note cell itself.
$wsPath.Rows.Item(7).RowHeight = 60
$wsPath.Range("D7").WrapText = $true

# New columns are now in use - size them like the author did.
$wsPath.Columns.Item(1).ColumnWidth = 55.3
$wsPath.Columns.Item(4).ColumnWidth = 71.65

# --- Active sheet / selection moved from hex_obj to hex_path -------------
$wsObj.Range("A5").Select()

$wsPath.Activate()
$wsPath.Range("D17").Select()
